$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Электроника"
$ws.Cells.Item(2, 2).Value = "—"

$ws.Cells.Item(3, 1).Value = "Телефоны"
$ws.Cells.Item(3, 2).Value = "Электроника"

$ws.Cells.Item(4, 1).Value = "Смартфоны"
$ws.Cells.Item(4, 2).Value = "Телефоны"

$ws.Cells.Item(5, 1).Value = "iPhone"
$ws.Cells.Item(5, 2).Value = "Смартфоны"

$ws.Cells.Item(6, 1).Value = "Samsung"
$ws.Cells.Item(6, 2).Value = "Смартфоны"

$ws.Cells.Item(7, 1).Value = "Electronics"
$ws.Cells.Item(7, 2).Value = "—"

$ws.Cells.Item(8, 1).Value = "Yers"
$ws.Cells.Item(8, 2).Value = "—"

$ws.Cells.Item(9, 1).Value = "sdu"
$ws.Cells.Item(9, 2).Value = "Yers"

$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "1"
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(10, 2).Value = "—"

$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "4"
$ws.Cells.Item(11, 1).ClearFormats()
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "1"
$ws.Cells.Item(11, 2).ClearFormats()

$ws.Cells.Item(12, 1).Value = "YersGay"
$ws.Cells.Item(12, 2).Value = "—"

$ws.Cells.Item(13, 1).Value = "Who"
$ws.Cells.Item(13, 2).Value = "YersGay"

$ws.Rows.Item(14).Delete()
